# Horarios actualizados Linea 141 - 874
# Refresh the scraped bus-schedule data across the three worksheets:
#   1) LP1912      - full schedule, oldest arrival drops off, new arrival appended
#   2) LP1912-215  - filtered schedule (stops containing "215"), one arrival drops off
#   3) 6203-6173   - no data rows, only the "last updated" timestamp changes

$wb = $excel.ActiveWorkbook

$newScrapTime = "03:24:43"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newScrapTime"
$ws1.Range("A3").Value = "Total filas: 6"

$sheet1Data = @(
    @($newScrapTime, "03:48", "14_ABASTO",      24,  "LP1912"),
    @($newScrapTime, "04:01", "81_EL PELIGRO",  37,  "LP1912"),
    @($newScrapTime, "04:45", "215A_EL PATO",   81,  "LP1912"),
    @($newScrapTime, "04:53", "11_ETCHEVERRY",  89,  "LP1912"),
    @($newScrapTime, "05:16", "17_ROMERO",      112, "LP1912"),
    @($newScrapTime, "05:22", "23_HERNANDEZ",   118, "LP1912")
)

$row = 6
foreach ($rec in $sheet1Data) {
    $ws1.Range("A$row").Value = $rec[0]
    $ws1.Range("B$row").Value = $rec[1]
    $ws1.Range("C$row").Value = $rec[2]
    $ws1.Range("D$row").Value = $rec[3]
    $ws1.Range("E$row").Value = $rec[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newScrapTime"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A6").Value = $newScrapTime
$ws2.Range("B6").Value = "04:45"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 81
$ws2.Range("E6").Value = "LP1912"

# Row 7 no longer exists in the refreshed data - remove it entirely.
$ws2.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newScrapTime"
